$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated K (formerly "Strike#") values for each row (column G), per updated
# std/mean calc. Map of row number -> new K value.
$updates = @{2=1; 3=0; 4=3; 5=2; 6=1; 7=1; 8=3; 9=0; 10=2; 11=1; 12=2; 13=2; 14=2; 15=2; 16=1; 17=2; 18=1; 19=2; 20=3; 21=1; 22=1; 23=1; 24=0; 25=1; 26=2; 27=2; 28=1; 29=2; 30=2; 32=2; 33=1; 34=2; 35=1; 36=2; 37=1; 38=0; 39=0; 40=1; 41=2; 42=0; 43=1; 44=1; 45=1; 46=2; 47=0; 48=2; 49=1; 50=2; 51=1; 52=2; 53=0; 54=0; 55=1; 56=2; 57=0; 58=1; 59=0; 60=1; 61=3; 62=2; 63=0; 64=1; 65=1; 66=1; 67=1; 68=1; 71=0; 72=1}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
